$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 102
$ws1.Range("G2").Value = 49
$ws1.Range("F3").Value = 12077
$ws1.Range("F4").Value = 43
$ws1.Range("F8").Value = 11968
$ws1.Range("F9").Value = 502
$ws1.Range("F10").Value = 1182
$ws1.Range("F11").Value = 112
$ws1.Range("F12").Value = 587
$ws1.Range("F13").Value = 1797
$ws1.Range("F14").Value = 5922
$ws1.Range("F16").Value = 3558
$ws1.Range("F17").Value = 201

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 578

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 578
$ws4.Range("F3").Value = 102
$ws4.Range("G3").Value = 49
$ws4.Range("F5").Value = 12077
$ws4.Range("F6").Value = 43
$ws4.Range("F11").Value = 11968
$ws4.Range("F12").Value = 502
$ws4.Range("F13").Value = 1182
$ws4.Range("F14").Value = 112
$ws4.Range("F15").Value = 587
$ws4.Range("F16").Value = 1797
$ws4.Range("F18").Value = 5922
$ws4.Range("F20").Value = 3558
$ws4.Range("F21").Value = 201
